$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.782.54"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "3.508.25"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'586.61"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'169.29"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.501.34"
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'0.575"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -5.07%  "
$ws.Range("D13").Value = "'46.96"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").Value = "4.074.17"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "'8.40"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("D17").Value = "'613.43"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  -9.04%  "
$ws.Range("D18").Value = "68.847.66"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "3.492.58"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "'17.40"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'11.09"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").Value = "'0.881"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -6.01%  "
$ws.Range("D24").Value = "'15.69"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("D25").Value = "'96.20"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -6.60%  "
$ws.Range("D29").Value = "'9.19"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -6.95%  "
$ws.Range("D30").Value = "'32.44"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  -6.34%  "
$ws.Range("E31").Value = "  -5.59%  "
$ws.Range("D32").Value = "'8.45"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -7.21%  "
$ws.Range("E33").Value = "  -5.38%  "
$ws.Range("E34").Value = "  -9.73%  "
$ws.Range("D35").Value = "'619.66"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +6.68%  "
$ws.Range("D36").Value = "'10.69"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("E37").Value = "  -4.90%  "
$ws.Range("D38").Value = "'57.06"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  -14.62%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'0.0443"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.379.97"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.134"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("D44").Value = "'0.325"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  -6.01%  "
$ws.Range("D45").Value = "'32.59"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  -5.40%  "
$ws.Range("D46").Value = "0.0₃0689"
$ws.Range("E46").Value = "  -5.70%  "
$ws.Range("E47").Value = "  -6.73%  "
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'133.32"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "'5.60"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +12.22%  "
